# Update the "Requisitos" (prerequisites) list on the active sheet:
#   - remove "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)"
#   - remove "LOQ4031 -  Química Geral I  (Requisito)"
#   - remove "LOQ4073 -  Química Geral II  (Requisito)"
#   - add "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)"
#   - add "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)"
#   - add "LOQ4247 -  Desenho Assistido por Computador  (Requisito)"
#
# The prerequisite rows live in columns B and C (both hold the same text),
# rows 24 through 49. The net change removes 3 entries and adds 3 entries,
# so the row range stays the same size; we just rewrite the list in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$requisitos = @(
    "LOB1003 -  Cálculo I  (Requisito)`n",
    "LOB1004 -  Cálculo II  (Requisito)`n",
    "LOB1006 -  Cálculo IV  (Requisito)`n",
    "LOB1011 -  Eletricidade Aplicada  (Requisito)`n",
    "LOB1012 -  Estatística  (Requisito)`n",
    "LOB1018 -  Física I  (Requisito)`n",
    "LOB1019 -  Física II  (Requisito)`n",
    "LOB1021 -  Física IV  (Requisito)`n",
    "LOB1024 -  Mecânica  (Requisito)`n",
    "LOB1036 -  Geometria Analítica  (Requisito)`n",
    "LOB1037 -  Àlgebra Linear  (Requisito)`n",
    "LOB1038 -  Física Experimental I  (Requisito)`n",
    "LOB1039 -  Física Experimental III  (Requisito)`n",
    "LOB1040 -  Laboratório de Eletricidade  (Requisito)`n",
    "LOB1041 -  Física Experimental II  (Requisito)`n",
    "LOB1042 -  Física Experimental IV  (Requisito)`n",
    "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)`n",
    "LOB1052 -  Cálculo III  (Requisito)`n",
    "LOB1053 -  Física III  (Requisito)`n",
    "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)`n",
    "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n",
    "LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito)`n",
    "LOQ4095 -  Química Geral Experimental  (Requisito)`n",
    "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)`n",
    "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)`n",
    "LOQ4247 -  Desenho Assistido por Computador  (Requisito)`n"
)

$startRow = 24
for ($i = 0; $i -lt $requisitos.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $requisitos[$i]
    $ws.Cells.Item($row, 3).Value = $requisitos[$i]
}
